# Fix for IAM test case
# Adds three new rows (125-127) of test case data to the "Test Cases" sheet,
# mirroring the style/formatting of the existing last row (124).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$steps = $wb.Worksheets.Item("Test Case Steps")

# ---- Row 125 ---------------------------------------------------------
$ws.Range("C124").Copy($ws.Range("A125"))
$ws.Range("A125").Value = "TestCase_B124"

$ws.Range("C124").Copy($ws.Range("B125"))
$ws.Range("B125").Value = "OPQA-583"

$ws.Range("C124").Copy($ws.Range("C125"))
$ws.Range("C125").Value = "Verify that record view page of an article gets displayed when user clicks on article title in ALL search results page"

$ws.Range("C124").Copy($ws.Range("D125"))
$ws.Range("D125").Value = "Y"

$ws.Range("C124").Copy($ws.Range("E125"))
$ws.Range("E125").ClearContents()

# ---- Row 126 ---------------------------------------------------------
$ws.Range("C124").Copy($ws.Range("A126"))
$ws.Range("A126").Value = "TestCase_B125"

$ws.Range("C124").Copy($ws.Range("B126"))
$ws.Range("B126").Value = "OPQA-585"

$ws.Range("C124").Copy($ws.Range("C126"))
$ws.Range("C126").Value = "Verify that record view page of an article gets displayed when user clicks on article title in ARTICLES search results page"

$ws.Range("C124").Copy($ws.Range("D126"))
$ws.Range("D126").Value = "Y"

$ws.Range("C124").Copy($ws.Range("E126"))
$ws.Range("E126").ClearContents()

# ---- Row 127 ---------------------------------------------------------
$ws.Range("C124").Copy($ws.Range("A127"))
$ws.Range("A127").Value = "TestCase_B126"

$ws.Range("C124").Copy($ws.Range("B127"))
$ws.Range("B127").Value = "OPQA-587"

# C127 uses the wrap-text style (s=4), so pull the format from the
# "Test Case Steps" sheet which already has cells using that style.
$steps.Range("A2").Copy($ws.Range("C127"))
$ws.Range("C127").Value = "Verify that following fields get displayed correctly for an article in record view page:`na)Title`nb)Publication name`nc)Publication date`nd)Publication volume`ne)Times cited count`nf)Cited references count`ng)Comments count`nh)Abstract`ng)DETAILS link"

$ws.Range("C124").Copy($ws.Range("D127"))
$ws.Range("D127").Value = "Y"

$ws.Range("C124").Copy($ws.Range("E127"))
$ws.Range("E127").ClearContents()

# The multi-line text in C127 would otherwise trigger an automatic
# "custom row height" recalculation; AutoFit() restores the row back to
# the sheet's standard (default) height with no explicit override.
$ws.Rows.Item(125).AutoFit()
$ws.Rows.Item(126).AutoFit()
$ws.Rows.Item(127).AutoFit()

# Update the visible selection/scroll position to match the edited area.
$ws.Range("C126").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 119
$excel.ActiveWindow.ScrollColumn = 1
